# Update DateBase/orders/Dang Nguyen 195_2026-2-9.xlsx
#
# - Orders sheet: row 61 gets new PackageID/FlowerName/Number values, and
#   five brand new rows (62-66) are appended with further flower lines,
#   the last of which (66) starts a new PackageID group (19).
# - Summary sheet: G2's running-totals string gets a new segment appended.

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell while forcing it to be stored as TEXT
# (matching t="str"/t="s" in the original file) instead of letting Excel
# auto-coerce numeric-looking strings (e.g. "18") into real numbers, and
# without losing precision on very long numeric-looking strings.
function Set-TextValue($range, [string]$value) {
    $escaped = $value -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

$ws = $wb.Worksheets.Item("Orders")

$orderRows = @(
    @{ Row = 61; A = "18"; C = "136_爱慕_adore_Rosa rugosa Thunb._20stems";        F = "20" },
    @{ Row = 62;           C = "203_佛罗伊德_Floyd_Rosa rugosa Thunb._20stems";     F = "4"  },
    @{ Row = 63;           C = "149_骄傲_Proud_Rosa rugosa Thunb._20stems";        F = "5"  },
    @{ Row = 64;           C = "208_紫霞仙子 _Nightingale_Rosa rugosa Thunb._20stems"; F = "5"  },
    @{ Row = 65;           C = "189_洛神_Mandala_Rosa rugosa Thunb._20stems";      F = "7"  },
    @{ Row = 66; A = "19"; C = "189_洛神_Mandala_Rosa rugosa Thunb._20stems";      F = "23" }
)

foreach ($r in $orderRows) {
    if ($r.ContainsKey("A")) {
        Set-TextValue $ws.Range("A" + $r.Row) $r.A
    }
    if ($r.ContainsKey("C")) {
        Set-TextValue $ws.Range("C" + $r.Row) $r.C
    }
    if ($r.ContainsKey("F")) {
        Set-TextValue $ws.Range("F" + $r.Row) $r.F
    }
}

$excel.CutCopyMode = $false

# Summary sheet: append the new encoded segment to the running total string.
$summary = $wb.Worksheets.Item("Summary")
$newG2 = "03014531467109145105338405302055501059570301001030738510121551542030101530312101051747.52015271018181020455723"
Set-TextValue $summary.Range("G2") $newG2

$excel.CutCopyMode = $false
